$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "João Rodrigues-Desenho Técnico"
$ws.Range("D2").Value = "José Ferreira-Tecnologia dos Materiais"
$ws.Range("F2").Value = "-"

$ws.Range("F3").Value = "-"

$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"

$ws.Range("E6").Value = "João Rodrigues-Desenho Técnico"
$ws.Range("F6").Value = "Andre Lucca-Circuitos Elétricos"

$ws.Range("C7").Value = "-"
$ws.Range("F7").Value = "Andre Lucca-Circuitos Elétricos"
